# Daily report auto-update: append new day row (row 26) to the report table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the last existing data row (25) down into the new row (26) first so
# the new row inherits the same per-column cell styles/number formats
# (date format on B, 2-decimal on G/H/I/N, integer on J/K) without minting
# new style entries in styles.xml.
$ws.Range("A25:Q25").Copy($ws.Range("A26:Q26"))

# --- Plain values -----------------------------------------------------
$ws.Range("A26").Value = 25
$ws.Range("B26").Value = 45743
$ws.Range("C26").Value = "Thu"
$ws.Range("D26").Value = 9163
$ws.Range("E26").Value = 9353
$ws.Range("L26").Value = 0
$ws.Range("O26").Value = 86736.1
$ws.Range("P26").Value = 87249.1

# --- Formulas (following the same pattern as the rows above) ----------
# Single-quoted so the '$' absolute-reference markers aren't treated as
# PowerShell variable expansions.
$ws.Range("F26").Formula = '=E26-D26'
$ws.Range("G26").Formula = '=(E26-$D$2)/A26'
$ws.Range("H26").Formula = '=(E26/D26-1)*100'
$ws.Range("I26").Formula = '=(POWER((E26/$D$3),1/A26)-1)*100'
$ws.Range("J26").Formula = '=J25*1.013'
$ws.Range("K26").Formula = '=E26-J26'
$ws.Range("M26").Formula = '=L26+E26'
$ws.Range("N26").Formula = '=E26/$D$2*100'
$ws.Range("Q26").Formula = '=P26/$O$2*100'

# Move the active selection to Q14, matching the saved session state.
$ws.Range("Q14").Select() | Out-Null
